$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.659.29"
$ws.Range("E2").Value = "  +5.45%  "

$ws.Range("D3").Value = "1.725.20"
$ws.Range("E3").Value = "  +3.69%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'225.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.73%  "

$ws.Range("D6").Value = "'0.5374"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.18%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.2700"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("D9").Value = "'0.06614"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.44%  "

$ws.Range("D10").Value = "'21.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.55%  "

$ws.Range("D11").Value = "'0.07756"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'4.656"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.738.11"
$ws.Range("E13").Value = "  +4.77%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.963.33"
$ws.Range("E14").Value = "  +3.73%  "

$ws.Range("D15").Value = "'0.5887"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.00%  "

$ws.Range("D16").Value = "0.0₅8300"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = "'68.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.41%  "

$ws.Range("D18").Value = "27.662.93"
$ws.Range("E18").Value = "  +5.49%  "

$ws.Range("D19").Value = "'224.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +16.10%  "

$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").Value = "'4.761"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "

$ws.Range("D22").Value = "'10.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.20%  "

$ws.Range("D23").Value = "'6.151"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.30%  "

$ws.Range("D24").Value = "'1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "'148.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "

$ws.Range("D26").Value = "'1.705"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.04%  "

$ws.Range("D27").Value = "'0.1235"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.34%  "

$ws.Range("D28").Value = "'7.433"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.26%  "

$ws.Range("D29").Value = "'16.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.56%  "

$ws.Range("D30").Value = "'0.05589"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").Value = "'1.303"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("D32").Value = "'3.599"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.82%  "

$ws.Range("D33").Value = "'3.475"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.25%  "

$ws.Range("D34").Value = "'1.667"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.69%  "

$ws.Range("D35").Value = "'0.9643"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "'2.827"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.43%  "

$ws.Range("D37").Value = "'2.450"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.91%  "

$ws.Range("D38").Value = "'0.5922"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("D39").Value = "'0.01656"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.16%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.8600"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.47%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.860"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.68%  "

$ws.Range("D42").Value = "1.062.41"
$ws.Range("E42").Value = "  +2.01%  "

$ws.Range("D43").Value = "'1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").Value = "'100.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("D45").Value = "1.869.44"
$ws.Range("E45").Value = "  +3.66%  "

$ws.Range("E46").Value = "  +10.34%  "

$ws.Range("D47").Value = "'59.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("D48").Value = "'8.222"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.30%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.4433"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").Value = "'1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.36%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05287"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
